$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table body: fill in the previously-blank "Resizable Array Bag" cells
# (C5:E5, C6:E6) with the same Big-O values/formatting already used by the
# "Linked Bag" columns (F:H), and correct the Worst-Case values that were
# wrong/missing.

# Row 5 (BestCase): match the formatting of F5 (blue fill) across C5:H5,
# all values are O(n+m).
$ws.Range("F5").Copy()
$ws.Range("C5:E5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C5").Value2 = "O(n+m)"
$ws.Range("D5").Value2 = "O(n+m)"
$ws.Range("E5").Value2 = "O(n+m)"
$ws.Range("F5").Value2 = "O(n+m)"
$ws.Range("G5").Value2 = "O(n+m)"
$ws.Range("H5").Value2 = "O(n+m)"

# Row 6 (WorstCase): match the formatting of F6 (blue fill) across C6:H6.
$ws.Range("F6").Copy()
$ws.Range("C6:E6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C6").Value2 = "O(n+m)"
$ws.Range("D6").Value2 = "O(m*n)"
$ws.Range("E6").Value2 = "O(m*n)"
$ws.Range("F6").Value2 = "O(n+m)"
$ws.Range("G6").Value2 = "O(m*n)"
$ws.Range("H6").Value2 = "O(m*n)"

# --- Comment on C3: replace the blank placeholder note with Dean Mah's
# write-up (mirrors the existing F3 comment from the other author).
$c3Text = @"
Dean Mah
Union Worst & Best Case - 
     Cloning the first bag will take O(n) since it has to copy every entry in the bag. It will take O(m) to add the entries from the second bag of size m to the clone bag. Therefore, it will take a total of O(n+m) since the bags can have different lengths.
Intersection Wost Case -
     The method will have to look at all entries in bag 2 for each of the entries in bag 1 making it O(n) from the first bag times the O(m) of the second bag therefore resulting in O(m*n) as it copies these over to a clone bag.
Intersection Best Case - 
     The method will only have to look at the first entry in bag 2 for each entry in bag 1 making the best case O(n+m) because it only has to look at each entry once and takes O(n) to clone them to the cloned return bag.
Difference Worst Case - 
     The method will take O(n) and O(m) to clone the bags. Then it will take the clone of the socond bag and look at each entry and remove it from the clone of the first bag if it contains the same entry before also removing it from the second bag. This process takes O(m*n) because it will have to look through every entry in the first bag looking to see if any are the same as the entry in the second bag which repeats for every entry in bag 2.
Difference Best Case - 
     If each entry in bag 2 matches the first entry of bag 1 it will only go over each entry in each bag once resulting in O(m+n) since it takes the time from going over the entries in bag 2 once O(m) and the time to clone the entries from the bags O(n).
"@
[void]$ws.Range("C3").Comment.Text($c3Text)

# --- Restore the active selection to A20, matching where the author left
# off editing.
[void]$ws.Range("A20").Select()
